$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '80.769.89'
$ws.Range("E2").Value = '  +5.76%  '

$ws.Range("D3").Value = '3.238.16'
$ws.Range("E3").Value = '  +6.27%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.30'
$ws.Range("E5").Value = '  +7.77%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '639.70'
$ws.Range("E6").Value = '  +3.39%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.292'
$ws.Range("E7").Value = '  +39.73%  '

$ws.Range("E8").Value = '  -0.08%  '

$ws.Range("E9").Value = '  +10.36%  '

$ws.Range("D10").Value = '3.238.42'
$ws.Range("E10").Value = '  +6.40%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.615'
$ws.Range("E11").Value = '  +39.51%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000272'
$ws.Range("E12").Value = '  +41.31%  '

$ws.Range("E13").Value = '  +3.60%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.50'
$ws.Range("E14").Value = '  +5.30%  '

$ws.Range("D15").Value = '3.834.55'
$ws.Range("E15").Value = '  +6.38%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '33.06'
$ws.Range("E16").Value = '  +14.21%  '

$ws.Range("D17").Value = '80.543.83'
$ws.Range("E17").Value = '  +5.63%  '

$ws.Range("D18").Value = '3.232.42'
$ws.Range("E18").Value = '  +6.57%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.76'
$ws.Range("E19").Value = '  +8.88%  '

$ws.Range("E20").Value = '  +25.74%  '

$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '451.66'
$ws.Range("E21").Value = '  +18.46%  '

$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.44'
$ws.Range("E22").Value = '  +5.75%  '

$ws.Range("E23").Value = '  +22.94%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.87'
$ws.Range("E24").Value = '  +12.36%  '

$ws.Range("D25").Value = '3.401.81'
$ws.Range("E25").Value = '  +7.04%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '77.98'
$ws.Range("E26").Value = '  +7.67%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.10'
$ws.Range("E27").Value = '  +12.80%  '

$ws.Range("E28").Value = '  +19.54%  '

$ws.Range("E29").Value = '  +0.08%  '

$ws.Range("E30").Value = '  +13.42%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  +0.28%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '569.78'
$ws.Range("E32").Value = '  +14.77%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.53'

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.159'
$ws.Range("E34").Value = '  +26.89%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.06'
$ws.Range("E35").Value = '  +7.58%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '23.88'
$ws.Range("E36").Value = '  +15.57%  '

$ws.Range("E37").Value = '  +21.46%  '

$ws.Range("E38").Value = '  +11.32%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.998'
$ws.Range("E39").Value = '  -0.16%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.93'
$ws.Range("E40").Value = '  +15.96%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '165.09'
$ws.Range("E41").Value = '  +1.61%  '

$ws.Range("E42").Value = '  +1.48%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '194.00'
$ws.Range("E43").Value = '  +1.03%  '

$ws.Range("E44").Value = '  +0.08%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.84'
$ws.Range("E45").Value = '  +16.56%  '

$ws.Range("E46").Value = '  +13.56%  '

$ws.Range("E47").Value = '  +9.44%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.809'
$ws.Range("E48").Value = '  +2.99%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '43.74'
$ws.Range("E49").Value = '  +6.18%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.39'
$ws.Range("E50").Value = '  +13.40%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.651'
$ws.Range("E51").Value = '  +9.62%  '
